$wb = $excel.ActiveWorkbook

# ---- Sheet: 单项选择题 (Single-choice) ----
$ws1 = $wb.Worksheets.Item("单项选择题")
$ws1.Range("A2").Value = "应用Johnson法则的流水作业调度采用的算法是"
$ws1.Range("B2").Value = "贪心算法"
$ws1.Range("C2").Value = "分支限界法"
$ws1.Range("D2").Value = "分治算法"
$ws1.Range("E2").Value = "动态规划法"
$ws1.Range("F2").Value = "D"

$ws1.Range("A3").Value = "动态规划算法的基本要素为"
$ws1.Range("B3").Value = "最优子结构性质与贪心选择性质"
$ws1.Range("C3").Value = "重叠子问题性质与贪心选择性质"
$ws1.Range("D3").Value = "最优子结构性质与重叠子问题性质"
$ws1.Range("E3").Value = "预排序与递归调用"
$ws1.Range("F3").Value = "C"

$ws1.Range("A4").Value = "能采用贪心算法求最优解的问题，一般具有的重要性质为："
$ws1.Range("B4").Value = "最优子结构性质与贪心选择性质"
$ws1.Range("C4").Value = "重叠子问题性质与贪心选择性质 "
$ws1.Range("D4").Value = "最优子结构性质与重叠子问题性质 "
$ws1.Range("E4").Value = "预排序与递归调用 "
$ws1.Range("F4").Value = "A"

$ws1.Range("A5").Value = "回溯法在问题的解空间树中，按（）策略，从根结点出发搜索解空间树。"
$ws1.Range("B5").Value = "广度优先"
$ws1.Range("C5").Value = "活结点优先"
$ws1.Range("D5").Value = "扩展节点优先"
$ws1.Range("E5").Value = "深度优先"
$ws1.Range("F5").Value = "D"

$ws1.Range("A6").Value = "分支限界法在问题的解空间树中，按（）策略，从根结点出发搜索解空间树。"
$ws1.Range("B6").Value = "广度优先"
$ws1.Range("C6").Value = "活结点优先"
$ws1.Range("D6").Value = "扩展节点优先"
$ws1.Range("E6").Value = "深度优先"
$ws1.Range("F6").Value = "A"

$ws1.Range("A7").Value = " 回溯法的效率不依赖于以下哪一个因素？"
$ws1.Range("B7").Value = "产生x[k]的时间；"
$ws1.Range("C7").Value = "满足显约束的x[k]值的个数； "
$ws1.Range("D7").Value = "问题的解空间的形式；  "
$ws1.Range("E7").Value = "计算上界函数bound的时间；"
$ws1.Range("F7").Value = "C"

$ws1.Range("A8").Value = "常见的两种分支限界法为"
$ws1.Range("B8").Value = "广度优先分支限界法与深度优先分支限界法；"
$ws1.Range("C8").Value = "队列式（FIFO）分支限界法与堆栈式分支限界法； "
$ws1.Range("D8").Value = "排列树法与子集树法； "
$ws1.Range("E8").Value = "队列式（FIFO）分支限界法与优先队列式分支限界法； "
$ws1.Range("F8").Value = "D"

$ws1.Range("A9").Value = "k带图灵机的空间复杂性S(n)是指"
$ws1.Range("B9").Value = "k带图灵机处理所有长度为n的输入时，在某条带上所使用过的最大方格数。"
$ws1.Range("C9").Value = "k带图灵机处理所有长度为n的输入时，在k条带上所使用过的方格数的总和"
$ws1.Range("D9").Value = "k带图灵机处理所有长度为n的输入时，在k条带上所使用过的平均方格数。 "
$ws1.Range("E9").Value = "k带图灵机处理所有长度为n的输入时，在某条带上所使用过的最小方格数。"
$ws1.Range("F9").Value = "B"

$ws1.Range("A10").Value = "NP类语言在图灵机下的定义为"
$ws1.Range("B10").Value = "NP={L|L是一个能在非多项式时间内被一台NDTM所接受的语言}"
$ws1.Range("C10").Value = "NP={L|L是一个能在多项式时间内被一台NDTM所接受的语言}"
$ws1.Range("D10").Value = "NP={L|L是一个能在多项式时间内被一台DTM所接受的语言}"
$ws1.Range("E10").Value = "NP={L|L是一个能在多项式时间内被一台NDTM所接受的语言}；"
$ws1.Range("F10").Value = "D"

$ws1.Range("A11").Value = "求最短路径可以使用什么算法"
$ws1.Range("B11").Value = "Prim算法"
$ws1.Range("C11").Value = "Dijkstra算法"
$ws1.Range("D11").Value = "Kruskal算法"
$ws1.Range("E11").Value = "冒泡排序算法"
$ws1.Range("F11").Value = "B"

# ---- Sheet: 多项选择题 (Multi-choice) ----
$ws2 = $wb.Worksheets.Item("多项选择题")
$ws2.Range("A2").Value = "算法设计多项选择题1"
$ws2.Range("F2").Value = "ABC"
$ws2.Range("A3").Value = "算法设计多项选择题2"
$ws2.Range("F3").Value = "BCD"
$ws2.Range("A4").Value = "算法设计多项选择题3"
$ws2.Range("F4").Value = "ABC"
$ws2.Range("A5").Value = "算法设计多项选择题4"
$ws2.Range("F5").Value = "BCD"

# ---- Sheet: 判断题 (True/False) ----
$ws3 = $wb.Worksheets.Item("判断题")
$ws3.Range("A2").Value = "算法设计判断题1"
$ws3.Range("B2").Value = "T"
$ws3.Range("A3").Value = "算法设计判断题2"
$ws3.Range("B3").Value = "F"
$ws3.Range("A4").Value = "算法设计判断题3"
$ws3.Range("B4").Value = "T"
$ws3.Range("A5").Value = "算法设计判断题4"
$ws3.Range("B5").Value = "F"
$ws3.Range("A6").Value = "算法设计判断题5"
$ws3.Range("B6").Value = "T"

# ---- Sheet: 主观题 (Subjective) ----
$ws4 = $wb.Worksheets.Item("主观题")
$ws4.Range("A2").Value = "试述回溯法的基本思想及用回溯法解题的步骤。"
$ws4.Range("B2").Value = "回溯算法答案。"
$ws4.Range("A3").Value = "简述归并排序算法和快速排序算法的分治方法。"
$ws4.Range("B3").Value = "归并与快排答案。"

# ---- Selections / active sheet ----
$ws2.Range("A6").Select()
$ws3.Range("C15").Select()
$ws4.Range("B3").Select()
$ws1.Range("F11").Select()